# Colony 4 revision: add a new "Snow Storm" weather entry to the language table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 89 (shifts existing rows 89.. down by one)
$ws.Rows.Item(89).Insert()

# Populate the new row with the new localization key/value pair
# (the inserted row already inherits the "wrap text" style used by column B
# from the row above, matching the rest of the table)
$ws.Cells.Item(89, 1).Value = "weather_snow_storm"
$ws.Cells.Item(89, 2).Value = "Snow Storm"

# Reflect the editor's final selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B89").Select() | Out-Null
